$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 already carries the exact style pattern (date / wrapped text x5 /
# mood) that the two new diary rows need, so clone its formatting down into
# rows 15 and 16 before writing the new content.
$ws.Range("A12:G12").Copy()
$ws.Range("A15:G16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 15: sixth diary entry (2/13/2020) ---
$ws.Range("A15").Value = "2/13/2020"
$ws.Range("B15").Value = "17:00-19:00"
$ws.Range("C15").Value = "Guowei Li, Dongxin Xiang"
$ws.Range("D15").Value = "Another three Key expert practices & Big picture of development"
$ws.Range("E15").Value = "Gained knowledges about the three Key expert practices & Big picture of development"
$ws.Range("F15").Value = "I had a idea how to find info about an open source project and how people contribute to it."
$ws.Range("G15").Value = "Tired because of the midterm, since it was longder than I expected. And I did not do well in the next day's job interview since I was tired.  I wish I could handle the situation better. Glad to learn something useful."

# --- Row 16: seventh diary entry (2/18-2/20/2020) ---
$ws.Range("A16").Value = "2/18/2020, 2/19/2020, 2/20/2020"
$ws.Range("B16").Value = "14:30-17:00, 9:00-15:00 (2/20/2020)"
$ws.Range("C16").Value = "Guowei Li, Dongxin Xiang"
$ws.Range("D16").Value = "Modify hw2 and complete hw3"
$ws.Range("E16").Value = 'Figure out how to improve hw2 with talking with Kaj. And we learned how to find "big picture" information about an open source project.'
$ws.Range("F16").Value = "Our group tried to figure out how to do the hw3. The requirement is simple, but we only learned the definations in class. So we need to figure out how and where to find related information.   And through teamwork, we get it done."
$ws.Range("G16").Value = "Pressured. I am still looking for an internship, and often have online assessments and interview during this quarter. As for me, it is not easy for these homework, since it took a lot of time. I have tried to make it balenced as possible as I can since I know both are important. But I feel more and more pressured recently. "

# Row heights grow to fit the newly-entered multi-line text.
$ws.Rows("15:15").RowHeight = 86.4
$ws.Rows("16:16").RowHeight = 159

# The formerly-blank row 17 loses its stray leftover formatting on A:C so it
# is uniform with D:G again.
$ws.Range("D17:F17").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Reflect the author's final view/selection state on the new row.
$ws.Range("G16").Select() | Out-Null
